$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# D1: company_name -> job_title
$ws.Range("D1").Value = "job_title"
# New E1: company (shifts old E/F/G content right)
$ws.Range("E1").Value = "company"
# F1: keep custom_message (was E1, now moved to F1)
$ws.Range("F1").Value = "custom_message"
# G1: sender_title -> email_sent_date
$ws.Range("G1").Value = "email_sent_date"
# New H1: status
$ws.Range("H1").Value = "status"

# Copy the header style (bold, border, centered) from an existing header
# cell onto the two newly added header cells so they match the rest.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# --- Data row (row 2) ---
$ws.Range("A2").Value = "your.email@example.com"
$ws.Range("B2").Value = "Your"
$ws.Range("C2").Value = "Name"
$ws.Range("D2").Value = "Your Title"
$ws.Range("E2").Value = "Your Company"
$ws.Range("F2").Value = "Your custom message here."

# G2 (email_sent_date) is left blank for this template row, but the cell
# itself is still present (not deleted) -- touch a formatting property
# that matches the existing default so the cell stays materialized
# without picking up a new style.
$ws.Range("G2").Value = ""
$ws.Range("G2").Font.Bold = $false

$ws.Range("H2").Value = "failed"
